$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 311, shifting existing rows 311-377 down to 312-378.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new record.
$ws.Range("A311").Value = 6
$ws.Range("B311").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C311").Value = "Metropolitana"
$ws.Range("D311").Value = 44511
$ws.Range("E311").Value = 13
$ws.Range("F311").Value = 100112044
$ws.Range("G311").Value = "Perejil"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 150
$ws.Range("K311").Value = 12000
$ws.Range("L311").Value = 13000
$ws.Range("M311").Value = 12400
$ws.Range("N311").Value = "$/docena de atados"
$ws.Range("O311").Value = "Región Metropolitana"
$ws.Range("P311").Value = 4133
$ws.Range("Q311").Value = 3
$ws.Range("R311").Value = "Hortaliza"
